$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stale "_GoBack" bookmark that currently sits right after
#    the inline drawing near the end of the document.  Word re-points
#    "_GoBack" at wherever editing last happened, so it will be
#    re-created below, at the surname fix.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Fix the misspelled surname "Трембыцький" -> "Трембицкий".
#    The edit point (right after "Трембиц") is where the new
#    "_GoBack" bookmark belongs, so split the run there.
# ------------------------------------------------------------------
$full = $d.Content
$null = $full.Find.Execute("Трембыцький")
$start = $full.Start

# "Трембиц" - first part of the fixed surname.
$r1 = $d.Range($start, $start + 7)
$r1.Text = "Трембиц"

# Drop "_GoBack" right after the edited prefix.
$mid = $start + 7
$bmRange = $d.Range($mid, $mid)
$d.Bookmarks.Add("_GoBack", $bmRange)

# "кий" - remaining part of the fixed surname (replaces old "ький").
$r2 = $d.Range($mid, $start + 11)
$r2.Text = "кий"

# ------------------------------------------------------------------
# 3) The text engine coalesces adjacent same-formatted runs in the
#    touched paragraph as part of the edit above, which also merges
#    the unrelated " Н.В" / "." runs later in the same paragraph into
#    " Н.В.". That merge is not part of the intended change, so split
#    those two runs back apart using a transient bookmark as a split
#    barrier, then remove the helper bookmark (bookmark removal alone
#    does not re-trigger the coalescing pass).
# ------------------------------------------------------------------
$tail = $d.Content
$gotTail = $tail.Find.Execute(" Н.В.")
if ($gotTail) {
    $dotStart = $tail.End - 1

    $barrier = $d.Range($dotStart, $dotStart)
    $d.Bookmarks.Add("ZZ_TempSplit", $barrier)
    $d.Bookmarks("ZZ_TempSplit").Delete()
}
